$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for the 380 kV case (Case_2_45), rows 2-25, columns B-N (excl. G/H)
$updates = @{
    "B2" = 1.02
    "C2" = 1.048516231170885
    "D2" = 1.054504249412253
    "E2" = 0.992614727750844
    "F2" = 1.063793041996202
    "I2" = 1.049377877516245
    "J2" = 1.053559568181004
    "K2" = 1.057247277639582
    "L2" = 0.9955398523335997
    "M2" = 1.066510721681059
    "N2" = 1.021675942880561
    "B3" = 1.02
    "C3" = 1.049470934283593
    "D3" = 1.055254880335033
    "E3" = 0.9936372048519299
    "F3" = 1.064682443636207
    "I3" = 1.049672388602259
    "J3" = 1.05416315428555
    "K3" = 1.057811324965967
    "L3" = 0.9963617723202687
    "M3" = 1.067215016665882
    "N3" = 1.021878645219127
    "B4" = 1.02
    "C4" = 1.050088851454443
    "D4" = 1.055740722588432
    "E4" = 0.9942998659930998
    "F4" = 1.065258404372407
    "I4" = 1.049861751975445
    "J4" = 1.054553239534588
    "K4" = 1.058175769136097
    "L4" = 0.9968940712668347
    "M4" = 1.06767056180101
    "N4" = 1.022009590913069
    "B5" = 1.02
    "C5" = 1.050348661478598
    "D5" = 1.055945001556694
    "E5" = 0.994578699834602
    "F5" = 1.065500646798244
    "I5" = 1.049941071346421
    "J5" = 1.054717116795602
    "K5" = 1.058328853232621
    "L5" = 0.9971179600053012
    "M5" = 1.067862028692796
    "N5" = 1.022064588389432
    "B6" = 1.02
    "C6" = 1.050392286910435
    "D6" = 1.055979302671706
    "E6" = 0.994625531979634
    "F6" = 1.065541326696668
    "I6" = 1.049954372456353
    "J6" = 1.054744625767131
    "K6" = 1.058354549182261
    "L6" = 0.9971555583673455
    "M6" = 1.067894174200705
    "N6" = 1.022073819643717
    "B7" = 1.02
    "C7" = 1.050092322901721
    "D7" = 1.055743452054103
    "E7" = 0.994303590798249
    "F7" = 1.065261640802616
    "I7" = 1.049862812980939
    "J7" = 1.054555429721763
    "K7" = 1.058177815158194
    "L7" = 0.9968970624462089
    "M7" = 1.067673120366238
    "N7" = 1.022010325997128
    "B8" = 1.02
    "C8" = 1.048838843885407
    "D8" = 1.054757900235131
    "E8" = 0.9929600610674297
    "F8" = 1.064093523802033
    "I8" = 1.049477658109536
    "J8" = 1.053763650877621
    "K8" = 1.057438009763721
    "L8" = 0.9958175282591056
    "M8" = 1.06674877830924
    "N8" = 1.021744491675213
    "B9" = 1.02
    "C9" = 1.046631315627512
    "D9" = 1.053022308199548
    "E9" = 0.9906006454969559
    "F9" = 1.062038716486958
    "I9" = 1.048789759802799
    "J9" = 1.05236482949684
    "K9" = 1.056130338910597
    "L9" = 0.9939188001724441
    "M9" = 1.065118628480412
    "N9" = 1.021274415186583
    "B10" = 1.02
    "C10" = 1.045160521617357
    "D10" = 1.051866038157421
    "E10" = 0.989033133672735
    "F10" = 1.060671311819824
    "I10" = 1.048325000910657
    "J10" = 1.05142990286524
    "K10" = 1.055255890820183
    "L10" = 0.9926553831429383
    "M10" = 1.064031018208086
    "N10" = 1.020959945010322
    "B11" = 1.02
    "C11" = 1.044523871492465
    "D11" = 1.051365562099962
    "E11" = 0.988355674866747
    "F11" = 1.060079810103435
    "I11" = 1.048122300114068
    "J11" = 1.051024514217649
    "K11" = 1.054876622271191
    "L11" = 0.9921088820399291
    "M11" = 1.063559882672446
    "N11" = 1.020823522170649
    "B12" = 1.02
    "C12" = 1.044287423954133
    "D12" = 1.051179693613923
    "E12" = 0.9881042295826724
    "F12" = 1.059860190322985
    "I12" = 1.048046789533899
    "J12" = 1.050873851490381
    "K12" = 1.054735651625772
    "L12" = 0.9919059725120875
    "M12" = 1.063384853849577
    "N12" = 1.020772810602294
    "B13" = 1.02
    "C13" = 1.044338141265833
    "D13" = 1.051219561632788
    "E13" = 0.9881581567098651
    "F13" = 1.059907295407411
    "I13" = 1.04806299668926
    "J13" = 1.050906172916805
    "K13" = 1.054765894523236
    "L13" = 0.9919494934313052
    "M13" = 1.06342239938363
    "N13" = 1.020783690123166
    "B14" = 1.02
    "C14" = 1.044504325995315
    "D14" = 1.051350197522322
    "E14" = 0.9883348863814464
    "F14" = 1.060061654415416
    "I14" = 1.048116062840321
    "J14" = 1.051012062085618
    "K14" = 1.05486497149509
    "L14" = 0.9920921077337197
    "M14" = 1.063545415303498
    "N14" = 1.020819331111172
    "B15" = 1.02
    "C15" = 1.044606722238537
    "D15" = 1.051430690693288
    "E15" = 0.9884438009545853
    "F15" = 1.0601567720931
    "I15" = 1.048148729709933
    "J15" = 1.051077292875354
    "K15" = 1.054926003730672
    "L15" = 0.9921799884222134
    "M15" = 1.063621205774309
    "N15" = 1.020841285666555
    "B16" = 1.02
    "C16" = 1.045202778536858
    "D16" = 1.051899257337632
    "E16" = 0.9890781214508737
    "F16" = 1.060710580441477
    "I16" = 1.048338422832271
    "J16" = 1.051456795450595
    "K16" = 1.055281048483135
    "L16" = 0.9926916645766087
    "M16" = 1.064062281930575
    "N16" = 1.020968993593009
    "B17" = 1.02
    "C17" = 1.045576726461894
    "D17" = 1.052193230188446
    "E17" = 0.989476357848556
    "F17" = 1.061058129484109
    "I17" = 1.048457022510265
    "J17" = 1.051694698148279
    "K17" = 1.055503591400877
    "L17" = 0.9930127773699352
    "M17" = 1.064338906288127
    "N17" = 1.021049033284879
    "B18" = 1.02
    "C18" = 1.045794864507878
    "D18" = 1.052364718465642
    "E18" = 0.9897087662937556
    "F18" = 1.061260906152207
    "I18" = 1.048526059074579
    "J18" = 1.05183340876158
    "K18" = 1.055633336397007
    "L18" = 0.9932001317071769
    "M18" = 1.064500237902724
    "N18" = 1.021095694442807
    "B19" = 1.02
    "C19" = 1.045869247411269
    "D19" = 1.052423194707079
    "E19" = 0.9897880325774034
    "F19" = 1.061330057422325
    "I19" = 1.048549574893822
    "J19" = 1.051880696314992
    "K19" = 1.055677565799461
    "L19" = 0.9932640239640975
    "M19" = 1.064555244615913
    "N19" = 1.021111600503719
    "B20" = 1.02
    "C20" = 1.045536603258728
    "D20" = 1.052161687722258
    "E20" = 0.9894336180360679
    "F20" = 1.061020834837259
    "I20" = 1.048444312429164
    "J20" = 1.051669179021322
    "K20" = 1.055479720906665
    "L20" = 0.9929783193494215
    "M20" = 1.064309229039442
    "N20" = 1.02104044833
    "B21" = 1.02
    "C21" = 1.044455387827935
    "D21" = 1.051311727655316
    "E21" = 0.9882828385668249
    "F21" = 1.060016197033295
    "I21" = 1.048100442209038
    "J21" = 1.050980882659517
    "K21" = 1.054835798369168
    "L21" = 0.9920501090198102
    "M21" = 1.063509190949828
    "N21" = 1.020808836776642
    "B22" = 1.02
    "C22" = 1.043775774064637
    "D22" = 1.050777501635151
    "E22" = 0.9875604150241495
    "F22" = 1.059385064338646
    "I22" = 1.047882973519066
    "J22" = 1.050547641616405
    "K22" = 1.054430398585132
    "L22" = 0.9914670000341481
    "M22" = 1.063006013075563
    "N22" = 1.020662993387984
    "B23" = 1.02
    "C23" = 1.044136032033518
    "D23" = 1.051060687759696
    "E23" = 0.9879432794643023
    "F23" = 1.059719589775516
    "I23" = 1.047998377436456
    "J23" = 1.050777356330286
    "K23" = 1.054645359630589
    "L23" = 0.991776070289318
    "M23" = 1.063272772257145
    "N23" = 1.020740328486849
    "B24" = 1.02
    "C24" = 1.045554733143653
    "D24" = 1.05217594034673
    "E24" = 0.9894529299347244
    "F24" = 1.061037686506214
    "I24" = 1.048450056002225
    "J24" = 1.051680710183232
    "K24" = 1.055490507142553
    "L24" = 0.9929938892766442
    "M24" = 1.064322638968497
    "N24" = 1.021044327577727
    "B25" = 1.02
    "C25" = 1.047201860624571
    "D25" = 1.053470865821461
    "E25" = 0.9912096547607049
    "F25" = 1.062569503792302
    "I25" = 1.048968685929499
    "J25" = 1.052726880997417
    "K25" = 1.056468876869545
    "L25" = 0.9944092447426414
    "M25" = 1.065540214235128
    "N25" = 1.021396133964533
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"